$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 4-9 with refined (higher-precision) values ---
$ws.Range("B4").Value = 30.0227352094
$ws.Range("C4").Value = 122.2853480397
$ws.Range("D4").Value = 53.8333557168
$ws.Range("E4").Value = 69.22777716829999
$ws.Range("F4").Value = 22.6503712232
$ws.Range("G4").Value = 92.284179673
$ws.Range("H4").Value = 85.5829243968
$ws.Range("I4").Value = 33.0044748827
$ws.Range("J4").Value = 93.969701066
$ws.Range("K4").Value = 45.7519694698
$ws.Range("L4").Value = 223.8241534809
$ws.Range("M4").Value = 114.5818639015
$ws.Range("N4").Value = 78.51659280299999

$ws.Range("B5").Value = 35.4844
$ws.Range("C5").Value = 122.3235
$ws.Range("D5").Value = 55.3456
$ws.Range("E5").Value = 71.54819999999999
$ws.Range("F5").Value = 20.9108
$ws.Range("G5").Value = 94.2088
$ws.Range("H5").Value = 88.6725
$ws.Range("I5").Value = 28.5493
$ws.Range("J5").Value = 96.44289999999999
$ws.Range("K5").Value = 49.7494
$ws.Range("L5").Value = 231.3544
$ws.Range("M5").Value = 123.6993
$ws.Range("N5").Value = 79.9863

$ws.Range("B6").Value = 37.5182163746
$ws.Range("C6").Value = 123.7832057705
$ws.Range("D6").Value = 56.9423449287
$ws.Range("E6").Value = 73.7403128895
$ws.Range("F6").Value = 20.8431267035
$ws.Range("G6").Value = 95.6745414136
$ws.Range("H6").Value = 90.67640584279999
$ws.Range("I6").Value = 29.1301718207
$ws.Range("J6").Value = 98.04909180769999
$ws.Range("K6").Value = 53.0654051485
$ws.Range("L6").Value = 235.3982590261
$ws.Range("M6").Value = 128.6098480043
$ws.Range("N6").Value = 80.7593158434

$ws.Range("B7").Value = 41.0094857109
$ws.Range("C7").Value = 121.2643063627
$ws.Range("D7").Value = 55.2390210092
$ws.Range("E7").Value = 79.08949080159999
$ws.Range("F7").Value = 19.4984755469
$ws.Range("G7").Value = 97.6853400359
$ws.Range("H7").Value = 97.1740735946
$ws.Range("I7").Value = 20.1683014697
$ws.Range("J7").Value = 100.9183410113
$ws.Range("K7").Value = 54.98788334
$ws.Range("L7").Value = 243.1340305041
$ws.Range("M7").Value = 142.179007666
$ws.Range("N7").Value = 73.14196024420001

$ws.Range("B8").Value = 43.1551025884
$ws.Range("C8").Value = 122.8435541582
$ws.Range("D8").Value = 55.731223961
$ws.Range("E8").Value = 81.65446012530001
$ws.Range("F8").Value = 18.7160459405
$ws.Range("G8").Value = 99.19638661899999
$ws.Range("H8").Value = 98.15951744669999
$ws.Range("I8").Value = 19.4974710411
$ws.Range("J8").Value = 102.5092004963
$ws.Range("K8").Value = 59.4214423884
$ws.Range("L8").Value = 247.3960254183
$ws.Range("M8").Value = 148.2813879798
$ws.Range("N8").Value = 72.1564365556

$ws.Range("B9").Value = 44.8925593184322
$ws.Range("C9").Value = 122.972410994478
$ws.Range("D9").Value = 56.4751483619934
$ws.Range("E9").Value = 82.5780507199576
$ws.Range("F9").Value = 18.232106721199
$ws.Range("G9").Value = 99.7292432454896
$ws.Range("H9").Value = 100.724599417361
$ws.Range("I9").Value = 19.3033416471746
$ws.Range("J9").Value = 103.12801307813
$ws.Range("K9").Value = 62.0147964533772
$ws.Range("L9").Value = 248.736470716752
$ws.Range("M9").Value = 149.555710847222
$ws.Range("N9").Value = 72.9436547616519

# --- Add new row 10 (2021年), formatted like the previous year rows ---
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "2021年"
$ws.Range("B10").Value = 50.1
$ws.Range("C10").Value = 120.3
$ws.Range("D10").Value = 55.4
$ws.Range("E10").Value = 82.3
$ws.Range("F10").Value = 18.2
$ws.Range("G10").Value = 100.5
$ws.Range("H10").Value = 98.09999999999999
$ws.Range("I10").Value = 12.7
$ws.Range("J10").Value = 104.2
$ws.Range("K10").Value = 68.8
$ws.Range("L10").Value = 253.6
$ws.Range("M10").Value = 161.7
$ws.Range("N10").Value = 63.2

# --- Add new row 11 (2022年); some metrics not yet reported, left blank ---
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "2022年"
$ws.Range("B11").Value = 51.4140502788285
$ws.Range("C11").Value = 120.605613311057
$ws.Range("F11").Value = 17.8730120235081
$ws.Range("G11").Value = 100.605341797751
$ws.Range("J11").Value = 104.358673401063
$ws.Range("L11").Value = 253.964768132876
$ws.Range("M11").Value = 163.497427518846
$ws.Range("N11").Value = 63.3923591148273
$excel.CutCopyMode = $false
